# Weekly update: insert a new daily-price record for
# "Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Arveja Verde"
# as row 194, pushing the existing rows 194:282 down to 195:283.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 194 (shifts old rows 194-282 down to 195-283,
# and Excel auto-extends the used range / dimension to R283).
$ws.Rows(194).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(194, 1).Value = 6
$ws.Cells.Item(194, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(194, 3).Value = "Metropolitana"
$ws.Cells.Item(194, 4).Value = 44917
$ws.Cells.Item(194, 5).Value = 13
$ws.Cells.Item(194, 6).Value = 100112022
$ws.Cells.Item(194, 7).Value = "Arveja Verde"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 400
$ws.Cells.Item(194, 11).Value = 25000
$ws.Cells.Item(194, 12).Value = 27000
$ws.Cells.Item(194, 13).Value = 25850
$ws.Cells.Item(194, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(194, 15).Value = "Carahue"
$ws.Cells.Item(194, 16).Value = 1034
$ws.Cells.Item(194, 17).Value = 25
$ws.Cells.Item(194, 18).Value = "Hortaliza"
